$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 2.888251324518907
$ws.Range("E2").Value = 16.71456242210409
$ws.Range("F2").Value = 22.91516405752202
$ws.Range("G2").Value = 3.554465621914876
$ws.Range("O2").Value = 19.41803292077184

$ws.Range("D3").Value = 2.868465467081504
$ws.Range("E3").Value = 15.74734985263153
$ws.Range("F3").Value = 22.18566614097683
$ws.Range("G3").Value = 3.558084161067147
$ws.Range("O3").Value = 18.9225144024005

$ws.Range("D4").Value = 2.856048993626708
$ws.Range("E4").Value = 15.12794543687679
$ws.Range("F4").Value = 21.73354054194758
$ws.Range("G4").Value = 3.560419138670807
$ws.Range("O4").Value = 18.61841796641222

$ws.Range("D5").Value = 2.850921696415791
$ws.Range("E5").Value = 14.86937712186839
$ws.Range("F5").Value = 21.54856223081966
$ws.Range("G5").Value = 3.561399240755086
$ws.Range("O5").Value = 18.49473964450639

$ws.Range("D6").Value = 2.850066233152149
$ws.Range("E6").Value = 14.82607902599314
$ws.Range("F6").Value = 21.51781232752258
$ws.Range("G6").Value = 3.56156371561648
$ws.Range("O6").Value = 18.47422376162875

$ws.Range("D7").Value = 2.855980118445161
$ws.Range("E7").Value = 15.12448283059674
$ws.Range("F7").Value = 21.73104838485161
$ws.Range("G7").Value = 3.560432240795915
$ws.Range("O7").Value = 18.61674873665769

$ws.Range("D8").Value = 2.88148468102177
$ws.Range("E8").Value = 16.38650937450282
$ws.Range("F8").Value = 22.66469596995112
$ws.Range("G8").Value = 3.555689875988607
$ws.Range("O8").Value = 19.24726691713166

$ws.Range("D9").Value = 2.929380739128488
$ws.Range("E9").Value = 18.80736662390714
$ws.Range("F9").Value = 24.44895728506584
$ws.Range("G9").Value = 3.547282679132358
$ws.Range("O9").Value = 20.47631796316273

$ws.Range("D10").Value = 2.963291022126727
$ws.Range("E10").Value = 20.49547908228555
$ws.Range("F10").Value = 25.71537767688067
$ws.Range("G10").Value = 3.541642386682016
$ws.Range("O10").Value = 21.36392112697865

$ws.Range("D11").Value = 2.978442842350294
$ws.Range("E11").Value = 21.22102878584216
$ws.Range("F11").Value = 26.27892897065387
$ws.Range("G11").Value = 3.539191319903616
$ws.Range("O11").Value = 21.76229304255136

$ws.Range("D12").Value = 2.984141199597955
$ws.Range("E12").Value = 21.48972210935707
$ws.Range("F12").Value = 26.4903152199851
$ws.Range("G12").Value = 3.538279535533489
$ws.Range("O12").Value = 21.91221204852719

$ws.Range("D13").Value = 2.982915702786961
$ws.Range("E13").Value = 21.432122964211
$ws.Range("F13").Value = 26.44488234123341
$ws.Range("G13").Value = 3.538475177899578
$ws.Range("O13").Value = 21.87996832303203

$ws.Range("D14").Value = 2.978912436747176
$ws.Range("E14").Value = 21.24325551487781
$ws.Range("F14").Value = 26.29636144028492
$ws.Range("G14").Value = 3.539115979214817
$ws.Range("O14").Value = 21.77464658914566

$ws.Range("D15").Value = 2.976455205107188
$ws.Range("E15").Value = 21.1267811182075
$ws.Range("F15").Value = 26.20511928788498
$ws.Range("G15").Value = 3.539510618405868
$ws.Range("O15").Value = 21.71000761895233

$ws.Range("D16").Value = 2.962295308117487
$ws.Range("E16").Value = 20.4472131168466
$ws.Range("F16").Value = 25.67827694917969
$ws.Range("G16").Value = 3.541804868354966
$ws.Range("O16").Value = 21.33776399663437

$ws.Range("D17").Value = 2.953538534236047
$ws.Range("E17").Value = 20.01949224721978
$ws.Range("F17").Value = 25.35170260145048
$ws.Range("G17").Value = 3.543241618569028
$ws.Range("O17").Value = 21.10790372692032

$ws.Range("D18").Value = 2.9484760591982
$ws.Range("E18").Value = 19.76948715586671
$ws.Range("F18").Value = 25.16269573567674
$ws.Range("G18").Value = 3.544078805245714
$ws.Range("O18").Value = 20.97519431035859

$ws.Range("D19").Value = 2.946757555924873
$ws.Range("E19").Value = 19.68415305554161
$ws.Range("F19").Value = 25.0985072080034
$ws.Range("G19").Value = 3.54436412175732
$ws.Range("O19").Value = 20.93018037350944

$ws.Range("D20").Value = 2.954473378242675
$ws.Range("E20").Value = 20.06543668049263
$ws.Range("F20").Value = 25.38658970316628
$ws.Range("G20").Value = 3.543087556526535
$ws.Range("O20").Value = 21.13242566986918

$ws.Range("D21").Value = 2.980089358902777
$ws.Range("E21").Value = 21.29889453732572
$ws.Range("F21").Value = 26.34004199921622
$ws.Range("G21").Value = 3.538927316677321
$ws.Range("O21").Value = 21.80560873600809

$ws.Range("D22").Value = 2.996601823452128
$ws.Range("E22").Value = 22.0697454502808
$ws.Range("F22").Value = 26.95133001253829
$ws.Range("G22").Value = 3.53630378881938
$ws.Range("O22").Value = 22.24005346421711

$ws.Range("D23").Value = 2.98780972986414
$ws.Range("E23").Value = 21.66154396560658
$ws.Range("F23").Value = 26.62622303281881
$ws.Range("G23").Value = 3.537695321746733
$ws.Range("O23").Value = 22.00873637834043

$ws.Range("D24").Value = 2.954050822766922
$ws.Range("E24").Value = 20.04467798395765
$ws.Range("F24").Value = 25.37082113682241
$ws.Range("G24").Value = 3.543157173143929
$ws.Range("O24").Value = 21.1213410305704

$ws.Range("D25").Value = 2.91664908718516
$ws.Range("E25").Value = 18.14787061189904
$ws.Range("F25").Value = 23.97300802016278
$ws.Range("G25").Value = 3.549462285454343
$ws.Range("O25").Value = 20.14578213644704
